# Horarios Línea 141 - update scrape (07:48:35 -> 07:55:46)
$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "LP1912"
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Header refresh
$ws1.Range("A2").Value = "Última actualización: 07:55:46"
$ws1.Range("A3").Value = "Total filas: 83"

# Rows 46/47 swap order (same scrape, re-sorted)
$ws1.Range("A46").Value = "07:12:53"
$ws1.Range("B46").Value = "07:59"
$ws1.Range("C46").Value = "23_HERNANDEZ"
$ws1.Range("D46").Value = 47
$ws1.Range("E46").Value = "LP1912"

$ws1.Range("A47").Value = "06:33:46"
$ws1.Range("B47").Value = "07:59"
$ws1.Range("C47").Value = "11_ETCHEVERRY"
$ws1.Range("D47").Value = 86
$ws1.Range("E47").Value = "LP1912"

# Insert a new row before the old row 81 (old rows 81-86 shift to 82-87)
$ws1.Rows(81).Insert()
$ws1.Range("A81").Value = "07:55:46"
$ws1.Range("B81").Value = "09:21"
$ws1.Range("C81").Value = "16_SANTA ANA"
$ws1.Range("D81").Value = 86
$ws1.Range("E81").Value = "LP1912"

# Append a new final row (88)
$ws1.Range("A88").Value = "07:55:46"
$ws1.Range("B88").Value = "09:51"
$ws1.Range("C88").Value = "15_ABASTO"
$ws1.Range("D88").Value = 116
$ws1.Range("E88").Value = "LP1912"

# -----------------------------------------------------------------
# Sheet "LP1912-215"
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 07:55:46"

# -----------------------------------------------------------------
# Sheet "6203-6173"
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 07:55:46"
$ws3.Range("A3").Value = "Total filas: 9"

# Insert a new row before the old row 11 (old rows 11-13 shift to 12-14)
$ws3.Rows(11).Insert()
$ws3.Range("A11").Value = "07:55:46"
$ws3.Range("B11").Value = "08:26"
$ws3.Range("C11").Value = "215C_LA PLATA"
$ws3.Range("D11").Value = 31
$ws3.Range("E11").Value = "L6203"
